# Goof progress on the irf functionality
#
# - "Priors" sheet: last user selection there moves to B5 (no longer the
#   active/visible tab).
# - "Iterations Calculator" sheet becomes the active tab, selection stays on
#   F9.
# - On "Iterations Calculator": the target "End Time" (F7) is bumped to a
#   later date, and the "iters / sec" cell (B11) is overwritten with a plain
#   typed-in number (3.4) instead of the computed formula.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Priors"
$ws2 = $wb.Worksheets.Item(2)   # "Iterations Calculator"

# Update the "Priors" sheet selection first (it stays the non-active tab).
[void]$ws1.Activate()
[void]$ws1.Range("B5").Select()

# New target end time: 2020-12-18 22:00:00 -> serial date 44183.916666666664
$ws2.Range("F7").Value = 44183.916666666664

# Overwrite the "iters / sec" formula with a plain constant the user typed.
$ws2.Range("B11").Value = 3.4

# Finish on the "Iterations Calculator" sheet, selection on F9, so it ends
# up as the saved active tab.
[void]$ws2.Activate()
[void]$ws2.Range("F9").Select()
